# Atualiza notas dos alunos
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the two students that are no longer present (Jordan, Lucas) -
# rows 4 and 5. Deleting them shifts Rosa (row 6) up to row 4.
$ws.Rows.Item(4).Delete() | Out-Null
$ws.Rows.Item(4).Delete() | Out-Null

# Insert a new grade column ("R5") before the Total column (currently F).
$ws.Columns.Item(6).Insert() | Out-Null

# New column header
$ws.Cells.Item(1, 6).Value = "R5"
$ws.Cells.Item(1, 6).Font.Bold = $true
$ws.Columns.Item(6).ColumnWidth = 5.1640625

# Fill in grade values
$ws.Cells.Item(2, 2).Value = 0
$ws.Cells.Item(2, 3).Value = 2

$ws.Cells.Item(3, 2).Value = 0
$ws.Cells.Item(3, 3).Value = 0

$ws.Cells.Item(4, 2).Value = 0
$ws.Cells.Item(4, 3).Value = 0

# Fix the Total formulas so they sum across the new R5 column too
$ws.Range("G2:G4").Formula = "=SUM(B2:F2)"

# Update selection to match the target state
$ws.Range("E5").Select() | Out-Null
